# Generate Report for Handback
# Updates the existing handed-back file's GUID/timestamps and appends a new
# handed-back file (a997d19f-...) as an additional row on every sheet.

$wb = $excel.ActiveWorkbook

$oldGuid = "e9de683e-e0e0-4cc4-9f69-9c4bce5cc8e8"
$newGuid = "365764c5-d128-40bc-9cee-edb6cb33f643"
$addGuid = "a997d19f-6a67-4018-8d32-d9177a7f1463"

$oldZhHash = "e6950e04337b36d644658fb292379d469697dae1"
$newZhHash = "1ce20cdbf5cf3afd932102f751cd6f9bea017291"
$addHash   = "68579c22ce767fa74bb61e44336f2d9a51bbe549"

# ---------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# Row 2: rename the handed-back file and refresh its generate date
$ws.Range("A2").Value = "$newGuid.md"
$ws.Range("B2").Value = "e2e\$newGuid.md"
$ws.Range("G2").Value = "2016-08-31 21:17:08"
$ws.Hyperlinks.Item(1).Address = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/128d41c1fac76b3253e92ed3493488cdfeb3552b/e2e/$newGuid.md"

# Row 3: new handed-back file
$ws.Range("A3").Value = "$addGuid.md"
$ws.Range("B3").Value = "e2e\$addGuid.md"
$ws.Range("C3").Value = ".md"
$ws.Range("E3").Value = "Handed back: in sync with en-US"
$ws.Range("F3").Value = "Handed back: in sync with en-US"
$ws.Range("G3").Value = "2016-08-31 21:17:08"
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/128d41c1fac76b3253e92ed3493488cdfeb3552b/e2e/$addGuid.md", "", "", "e2e\$addGuid.md") | Out-Null

# ---------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

# Row 2: rename the handed-back file and refresh handoff/handback dates
$ws.Range("A2").Value = "$newGuid.md"
$ws.Range("G2").Value = "$newGuid.$newZhHash.zh-cn.xlf"
$ws.Range("H2").Value = "2016-08-31 21:16:58"
$ws.Range("I2").Value = "$newGuid.md"
$ws.Range("J2").Value = "$newGuid.$newZhHash.zh-cn.xlf"
$ws.Range("K2").Value = "2016-08-31 21:17:29"
$ws.Hyperlinks.Item(1).Address = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/128d41c1fac76b3253e92ed3493488cdfeb3552b/e2e/$newGuid.md"
$ws.Hyperlinks.Item(2).Address = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/298c8abdcc7e6348385a1137f0be1f28bf07b655/e2e/$newGuid.md"

# Row 3: new handed-back file
$ws.Range("A3").Value = "$addGuid.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("D3").Value = "e2e"
$ws.Range("E3").Value = "ht"
$ws.Range("F3").Value = "True"
$ws.Range("G3").Value = "$addGuid.$addHash.zh-cn.xlf"
$ws.Range("H3").Value = "2016-08-31 21:16:58"
$ws.Range("I3").Value = "$addGuid.md"
$ws.Range("J3").Value = "$addGuid.$addHash.zh-cn.xlf"
$ws.Range("K3").Value = "2016-08-31 21:17:29"
$ws.Range("L3").Value = ""
$ws.Range("M3").Value = "True"
$ws.Range("N3").Value = ""
$ws.Range("O3").Value = "False"
$ws.Range("P3").Value = ""
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/128d41c1fac76b3253e92ed3493488cdfeb3552b/e2e/$addGuid.md", "", "", "$addGuid.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/298c8abdcc7e6348385a1137f0be1f28bf07b655/e2e/$addGuid.md", "", "", "$addGuid.md") | Out-Null

# ---------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

# Row 2: rename the handed-back file and refresh handoff/handback dates
$ws.Range("A2").Value = "$newGuid.md"
$ws.Range("G2").Value = "$newGuid.$newZhHash.de-de.xlf"
$ws.Range("H2").Value = "2016-08-31 21:17:08"
$ws.Range("I2").Value = "$newGuid.md"
$ws.Range("J2").Value = "$newGuid.$newZhHash.de-de.xlf"
$ws.Range("K2").Value = "2016-08-31 21:17:38"
$ws.Hyperlinks.Item(1).Address = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/128d41c1fac76b3253e92ed3493488cdfeb3552b/e2e/$newGuid.md"
$ws.Hyperlinks.Item(2).Address = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/da472420557c33366c9f2a3658b29d349060bb8d/e2e/$newGuid.md"

# Row 3: new handed-back file
$ws.Range("A3").Value = "$addGuid.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("D3").Value = "e2e"
$ws.Range("E3").Value = "ht"
$ws.Range("F3").Value = "True"
$ws.Range("G3").Value = "$addGuid.$addHash.de-de.xlf"
$ws.Range("H3").Value = "2016-08-31 21:17:08"
$ws.Range("I3").Value = "$addGuid.md"
$ws.Range("J3").Value = "$addGuid.$addHash.de-de.xlf"
$ws.Range("K3").Value = "2016-08-31 21:17:38"
$ws.Range("L3").Value = ""
$ws.Range("M3").Value = "True"
$ws.Range("N3").Value = ""
$ws.Range("O3").Value = "False"
$ws.Range("P3").Value = ""
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/128d41c1fac76b3253e92ed3493488cdfeb3552b/e2e/$addGuid.md", "", "", "$addGuid.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/da472420557c33366c9f2a3658b29d349060bb8d/e2e/$addGuid.md", "", "", "$addGuid.md") | Out-Null
